$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-13 23:48:24'
$ws.Range("I2").Value = '5.3 mm'
$ws.Range("E3").Value = '2026-02-13 23:48:26'
$ws.Range("I3").Value = '8.3 mm'
$ws.Range("L3").Value = '56.2 km/h - 262º 23:16 TU'
$ws.Range("E4").Value = '2026-02-13 23:48:29'
$ws.Range("J4").Value = '993.1 hPa'
$ws.Range("K4").Value = '3.1 MJ/m2'
$ws.Range("E5").Value = '2026-02-13 23:48:31'
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = '84%'
$ws.Range("I5").Value = '5.4 mm'
$ws.Range("N5").Value = '-4.5 °C 23:25 TU'
$ws.Range("O5").Value = '-2.9 °C'
$ws.Range("E6").Value = '2026-02-13 23:48:34'
$ws.Range("J6").Value = '993.1 hPa'
$ws.Range("E7").Value = '2026-02-13 23:48:36'
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = '71%'
$ws.Range("J7").Value = '993.4 hPa'
$ws.Range("O7").Value = '12.6 °C'
$ws.Range("E8").Value = '2026-02-13 23:48:39'
$ws.Range("J8").Value = '993.4 hPa'
$ws.Range("O8").Value = '8.9 °C'
$ws.Range("E9").Value = '2026-02-13 23:48:41'
$ws.Range("E10").Value = '2026-02-13 23:48:44'
$ws.Range("K10").Value = '2.4 MJ/m2'
$ws.Range("E11").Value = '2026-02-13 23:48:46'
$ws.Range("E12").Value = '2026-02-13 23:48:49'
$ws.Range("E13").Value = '2026-02-13 23:48:51'
$ws.Range("J13").Value = '996.1 hPa'
$ws.Range("E14").Value = '2026-02-13 23:48:54'
$ws.Range("O14").Value = '10.6 °C'
$ws.Range("E15").Value = '2026-02-13 23:48:56'
$ws.Range("I15").Value = '6.5 mm'
$ws.Range("E16").Value = '2026-02-13 23:48:59'
$ws.Range("I16").Value = '15.0 mm'
$ws.Range("O16").Value = '-4.1 °C'
$ws.Range("E17").Value = '2026-02-13 23:49:01'
$ws.Range("E18").Value = '2026-02-13 23:49:03'
$ws.Range("J18").Value = '993.3 hPa'
$ws.Range("E19").Value = '2026-02-13 23:49:06'
$ws.Range("E20").Value = '2026-02-13 23:49:08'
$ws.Range("I20").Value = '25.4 mm'
$ws.Range("E21").Value = '2026-02-13 23:49:11'
$ws.Range("J21").Value = '996.3 hPa'
$ws.Range("E22").Value = '2026-02-13 23:49:13'
$ws.Range("E23").Value = '2026-02-13 23:49:16'
$ws.Range("I23").Value = '14.5 mm'
$ws.Range("E24").Value = '2026-02-13 23:49:18'
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = '93%'
$ws.Range("J24").Value = '994.6 hPa'
$ws.Range("E25").Value = '2026-02-13 23:49:21'
$ws.Range("G25").Value = '117 cm'
$ws.Range("I25").Value = '10.9 mm'
$ws.Range("E26").Value = '2026-02-13 23:49:23'
$ws.Range("E27").Value = '2026-02-13 23:49:26'
$ws.Range("E28").Value = '2026-02-13 23:49:28'
$ws.Range("J28").Value = '993.5 hPa'
$ws.Range("E29").Value = '2026-02-13 23:49:31'
$ws.Range("E30").Value = '2026-02-13 23:49:33'
$ws.Range("J30").Value = '993.0 hPa'
$ws.Range("O30").Value = '9.4 °C'
$ws.Range("E31").Value = '2026-02-13 23:49:36'
$ws.Range("I31").Value = '6.0 mm'
$ws.Range("J31").Value = '992.0 hPa'
$ws.Range("O31").Value = '10.0 °C'
$ws.Range("E32").Value = '2026-02-13 23:49:39'
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = '92%'
$ws.Range("E33").Value = '2026-02-13 23:49:41'
$ws.Range("J33").Value = '995.1 hPa'
$ws.Range("E34").Value = '2026-02-13 23:49:44'
$ws.Range("I34").Value = '10.7 mm'
$ws.Range("E35").Value = '2026-02-13 23:49:46'
$ws.Range("J35").Value = '994.8 hPa'
$ws.Range("L35").Value = '94.3 km/h - 273º 23:15 TU'
$ws.Range("N35").Value = '3.5 °C 23:29 TU'
$ws.Range("E36").Value = '2026-02-13 23:49:49'
$ws.Range("I36").Value = '9.0 mm'
$ws.Range("J36").Value = '993.1 hPa'
$ws.Range("E37").Value = '2026-02-13 23:49:51'
$ws.Range("J37").Value = '995.0 hPa'
$ws.Range("E38").Value = '2026-02-13 23:49:53'
$ws.Range("N38").Value = '7.5 °C 23:26 TU'
$ws.Range("E39").Value = '2026-02-13 23:49:56'
$ws.Range("I39").Value = '20.3 mm'
$ws.Range("E40").Value = '2026-02-13 23:49:58'
$ws.Range("J40").Value = '996.8 hPa'
$ws.Range("E41").Value = '2026-02-13 23:50:01'
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = '73%'
$ws.Range("J41").Value = '993.9 hPa'
$ws.Range("E42").Value = '2026-02-13 23:50:03'
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = '84%'
$ws.Range("E43").Value = '2026-02-13 23:50:05'
$ws.Range("E44").Value = '2026-02-13 23:50:08'
$ws.Range("I44").Value = '12.1 mm'
$ws.Range("L44").Value = '78.8 km/h - 107º 23:15 TU'
$ws.Range("E45").Value = '2026-02-13 23:50:10'
$ws.Range("I45").Value = '4.2 mm'
$ws.Range("E46").Value = '2026-02-13 23:50:13'
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = '85%'
$ws.Range("J46").Value = '994.8 hPa'
$ws.Range("O46").Value = '9.4 °C'
